$d = $word.ActiveDocument
$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

# 1) Title paragraph: wrap "New_Sample_App" run with spellStart/spellEnd proofErr markers.
$p1 = $d.Paragraphs.Item(1)
$p1xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:jc w:val="center"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>New_Sample_App</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> Users</w:t></w:r></w:p>
'@
$p1.Range.InsertXML($p1xml)

# 2) "Name: Robrodjr" paragraph: split into "Name: " run + spellchecked "Robrodjr" run.
$p8 = $d.Paragraphs.Item(8)
$p8xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Name: </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Robrodjr</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>
'@
$p8.Range.InsertXML($p8xml)

# 3) After "password: supaidaman78" insert the new account records (and 8 extra blank
#    paragraphs before the pre-existing trailing blank paragraphs).
$p14 = $d.Paragraphs.Item(14)
$blockXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">New yahoo.com email for Sample Email for </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Heroku</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> production</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Mario Gomez</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>mar_gom_77</w:t></w:r><w:r><w:t>74</w:t></w:r><w:r><w:t>@yahoo.com</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">password: </w:t></w:r><w:r><w:t>nina0197</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>DOB 02/16/1970</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">New yahoo.com email for Sample Email for </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Heroku</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> production</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Rob Rod</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>rob_rod_jr@yahoo.com</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>password: mydob0216</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>DOB 02/16/1971</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">New yahoo.com email for Sample Email for </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Heroku</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> production</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:lastRenderedPageBreak/><w:t>John Parker</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>john_parker_1984@yahoo.com</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>password: davedob0204</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>DOB 02/16/1972</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
'@
$p14.Range.InsertXML($blockXml)

Write-Output "done; paragraph count = $($d.Paragraphs.Count)"
